$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = -12.726
$ws.Range("B3").Value = 6.406999999999999
$ws.Range("C5").Value = -12.589
$ws.Range("B14").Value = 6.645999999999999
$ws.Range("B16").Value = 5.44
$ws.Range("C16").Value = -12.245
$ws.Range("B21").Value = 6.179
$ws.Range("B23").Value = 6.444
$ws.Range("B25").Value = 6.103999999999999
